{"js": "// \"fix feeAmount in template\"\n// The template placeholder \"<<fee amount>>\" contained a stray space that\n// broke the merge field name. Fix it so the field reads \"<<feeAmount>>\".\nconst body = context.document.body;\n\nconst results = body.search(\"fee amount\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"feeAmount\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# \"fix feeAmount in template\"\n# The template placeholder \"<<fee amount>>\" contained a stray space that\n# broke the merge field name. Fix it so the field reads \"<<feeAmount>>\".\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"fee amount\"\n$find.Replacement.Text = \"feeAmount\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 0          # wdFindStop - do not wrap around the document\n\n# wdReplaceOne (1) - Find.Execute's own MatchCase/Forward/Wrap members above\n# drive the actual match; positional args just mirror them for clarity.\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 1)\n"}
